$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 64999.5
$ws.Range("J3").Value = 64999.5
$ws.Range("L3").Value = 64999.5
$ws.Range("N3").Value = -65227.5

$ws.Range("H4").Value = 17
$ws.Range("I4").Value = 17
$ws.Range("K4").Value = 17
$ws.Range("M4").Value = 97

$ws.Range("H28").Value = 433
$ws.Range("I28").Value = 422.57144
$ws.Range("J28").Value = 506
$ws.Range("K28").Value = 422.57144
$ws.Range("L28").Value = 506
$ws.Range("M28").Value = 62.42856
$ws.Range("N28").Value = -1476

$ws.Range("H29").Value = 2610.75
$ws.Range("J29").Value = 4022.6
$ws.Range("L29").Value = 12067.8
$ws.Range("N29").Value = -12629.8

$ws.Range("H41").Value = 2581.6
$ws.Range("I41").Value = 3262
$ws.Range("J41").Value = 994
$ws.Range("K41").Value = 3262
$ws.Range("L41").Value = 994
$ws.Range("M41").Value = -2822
$ws.Range("N41").Value = -1874

$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("N59").ClearContents()

$ws.Range("H102").Value = 64999.5
$ws.Range("J102").Value = 64999.5
$ws.Range("L102").Value = 64999.5
$ws.Range("N102").Value = -71489.5

$ws.Range("H129").Value = 6382.3335
$ws.Range("I129").Value = 637.2222
$ws.Range("K129").Value = 1911.6666
$ws.Range("M129").Value = 3088.3334

$ws.Range("H132").Value = 2378.889
$ws.Range("I132").Value = 2378.889
$ws.Range("K132").Value = 7136.667
$ws.Range("M132").Value = -4606.667

$ws.Range("I137").Value = 1967.1666
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 5901.4998
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3351.4998
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6214.4287
$ws.Range("I32").Value = 7200.2
$ws.Range("K32").Value = 7200.2
$ws.Range("M32").Value = -6913.2

$ws.Range("H110").Value = 1297.5
$ws.Range("I110").Value = 1297.5
$ws.Range("K110").Value = 1297.5
$ws.Range("M110").Value = 747.5

$ws.Range("H118").Value = 25000
$ws.Range("J118").Value = 25000
$ws.Range("L118").Value = 25000
$ws.Range("N118").Value = -28314

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 9999
$ws.Range("I105").Value = 9999
$ws.Range("K105").Value = 9999
$ws.Range("M105").Value = -8252

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 301.33334
$ws.Range("I6").Value = 301.33334
$ws.Range("K6").Value = 301.33334
$ws.Range("M6").Value = -188.33334

$ws.Range("H7").Value = 168.09091
$ws.Range("J7").Value = 200.2
$ws.Range("L7").Value = 200.2
$ws.Range("N7").Value = -426.2

$ws.Range("H16").Value = 1152.75
$ws.Range("I16").Value = 805.5
$ws.Range("K16").Value = 805.5
$ws.Range("M16").Value = -518.5

$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 500
$ws.Range("K22").Value = 500
$ws.Range("M22").Value = -150

$ws.Range("H31").Value = 2730.75
$ws.Range("J31").Value = 2631.6667
$ws.Range("L31").Value = 2631.6667
$ws.Range("N31").Value = -3221.6667

$ws.Range("H34").Value = 2730.75
$ws.Range("J34").Value = 2631.6667
$ws.Range("L34").Value = 2631.6667
$ws.Range("N34").Value = -3035.6667

$ws.Range("H88").Value = 35000
$ws.Range("J88").Value = 35000
$ws.Range("L88").Value = 35000
$ws.Range("N88").Value = -35812

$ws.Range("H91").Value = 35000
$ws.Range("J91").Value = 35000
$ws.Range("L91").Value = 35000
$ws.Range("N91").Value = -37808

$ws.Range("H108").Value = 75000
$ws.Range("J108").Value = 75000
$ws.Range("L108").Value = 75000
$ws.Range("N108").Value = -82680

$ws.Range("H113").Value = 1152.75
$ws.Range("I113").Value = 805.5
$ws.Range("K113").Value = 805.5
$ws.Range("M113").Value = 1364.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 47.25
$ws.Range("I6").Value = 39.714287
$ws.Range("K6").Value = 119.142861
$ws.Range("M6").Value = -6.142860999999996

$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").ClearContents()

$ws.Range("H69").Value = 1547
$ws.Range("J69").Value = 1547
$ws.Range("L69").Value = 4641
$ws.Range("N69").Value = -6263

$ws.Range("H72").Value = 1547
$ws.Range("J72").Value = 1547
$ws.Range("L72").Value = 13923
$ws.Range("N72").Value = -22035

$ws.Range("H113").Value = 762.8
$ws.Range("I113").Value = 604.6667
$ws.Range("K113").Value = 1814.0001
$ws.Range("M113").Value = 355.9999

$ws.Range("H131").Value = 1099.3
$ws.Range("I131").Value = 713.2857
$ws.Range("K131").Value = 2139.8571
$ws.Range("M131").Value = 2900.1429

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3000
$ws.Range("I7").Value = 3000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3000
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2888
$ws.Range("N7").ClearContents()

$ws.Range("H18").Value = 4647.6665
$ws.Range("J18").Value = 4647.6665
$ws.Range("L18").Value = 4647.6665
$ws.Range("N18").Value = -4991.6665

$ws.Range("H46").Value = 8999.6
$ws.Range("I46").Value = 4999.6665
$ws.Range("K46").Value = 4999.6665
$ws.Range("M46").Value = -4811.6665

$ws.Range("H55").Value = 278.42856
$ws.Range("I55").Value = 275.8
$ws.Range("K55").Value = 275.8
$ws.Range("M55").Value = -102.8

$ws.Range("H68").Value = 2002
$ws.Range("I68").Value = 2002
$ws.Range("K68").Value = 2002
$ws.Range("M68").Value = -1253

$ws.Range("H71").Value = 2002
$ws.Range("I71").Value = 2002
$ws.Range("K71").Value = 10010
$ws.Range("M71").Value = -6266

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -6530
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 10004
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 10004
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 10004
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -10340
